$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 24 ("福岡県保健環境研究所" 2012 entry), which shifts all rows below
# it up by one (old row 25 becomes new row 24, ..., old row 65 becomes new row 64).
$ws.Rows.Item(24).Delete()
